$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows: URL / Body columns now use "serviceDelivery" (camelCase)
# and JSON-ish bodies now use " : " instead of ", " as separator.
$ws.Cells.Item(2, 4).Value = '{"foo" : "bar"}'

$ws.Cells.Item(3, 3).Value = "/serviceDelivery/101"
$ws.Cells.Item(3, 4).Value = '{"serviceDelivery" : "101"}'

$ws.Cells.Item(4, 3).Value = "/serviceDelivery/102"
$ws.Cells.Item(4, 4).Value = '{"serviceDelivery" : "102"}'

$ws.Cells.Item(5, 3).Value = "/serviceDelivery/103"
$ws.Cells.Item(5, 4).Value = '{"serviceDelivery" : "103"}'

# Add new stub row for POST, TestNG test
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "POST"
$ws.Cells.Item(6, 3).Value = "/serviceDelivery"
$ws.Cells.Item(6, 5).Value = 200

# Selection moves to F6 after the edits
$ws.Range("F6").Select()

# Page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
